# Apply header-renaming edit: handle multiple primary keys
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row cell values (B1:D1); E1 ("df") stays the same.
$ws.Range("B1").Value = "Risk"
$ws.Range("C1").Value = "Cell"
$ws.Range("D1").Value = "Mask"

# Move the active selection to J11 (matches the saved sheet view state).
$ws.Range("J11").Select()
